$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the JSON-ish test inputs in column B (rows 2-9) to use proper
# double-quoted JSON syntax instead of the old single-quoted / unquoted keys.
$ws.Range("B2").Value = "API 1 - Add friend`nInput`n{   `"friends`":     `n[       `"andy@example.com`",       `"john@example.com`"     ] } "
$ws.Range("B3").Value = "API 1 - Add friend`nInput`n{   `"friends`":     `n[       `"andy`",       `"john`"     ] } "
$ws.Range("B4").Value = "API 2 - Get Friend List`nInput`n{   `"email`": `"andy@example.com`" `n} "
$ws.Range("B5").Value = "API 2 - Get Friend List`nInput`n{   `"email`": `"andy`"`n} "
$ws.Range("B6").Value = "API 1 - Add friend`nInput`n{   `"friends`":     `n[       `"andy@example.com`",       `"common@example.com`"] } "
$ws.Range("B7").Value = "API 1 - Add friend`nInput`n{   `"friends`":     `n[       `"john@example.com`",       `"common@example.com`"] } "
$ws.Range("B8").Value = "API 3 - Get Common Friends`nInput`n{   `"friends`":     `n[       `"andy@example.com`",       `"john@example.com`"     ] } "
$ws.Range("B9").Value = "API 3 - Get Common Friends`nInput`n{   `"friends`":     `n[       `"andy@example.com`",       `"kevin@example.com`"     ] } "

# Select cell B1 (matches the selection recorded in the saved sheet view).
$ws.Range("B1").Select() | Out-Null
